$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="42.438.18"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Formula = '="2.290.01"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Formula = '="301.69"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").Formula = '="95.69"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Formula = '="0.508"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Formula = '="34.46"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D11").Formula = '="18.97"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +3.59%  '
$ws.Range("D12").Formula = '="0.0782"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Formula = '="6.75"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("D16").Formula = '="2.281.15"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Formula = '="42.393.70"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").Formula = '="12.17"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -6.18%  '
$ws.Range("D20").Formula = '="0.0₃0886"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Formula = '="67.71"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").Formula = '="2.27"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +6.72%  '
$ws.Range("D24").Formula = '="235.31"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("E28").Value = '  +15.13%  '
$ws.Range("D29").Formula = '="165.73"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Formula = '="31.90"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -3.09%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").Formula = '="5.00"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("D34").Formula = '="17.54"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").Formula = '="4.42"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -7.00%  '
$ws.Range("E36").Value = '  +1.17%  '
$ws.Range("E37").Value = '  -2.68%  '
$ws.Range("D38").Formula = '="0.0997"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Formula = '="20.18"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +11.82%  '
$ws.Range("D43").Formula = '="1.964.07"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").Formula = '="10.39"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +4.06%  '
$ws.Range("D45").Formula = '="0.0279"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").Formula = '="2.515.58"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("D50").Formula = '="53.14"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").Formula = '="71.16"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -0.19%  '
